# Attendance export update: add a new date column (U) with "2025-11-22"
# header, mark every student absent ("❌") for that date, and bump the
# per-row "Total" (column S) by 1 to account for the new tracked date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell U1 - force text formatting so the date string isn't
# auto-converted into a serial date value, then copy the neighboring
# header cell's formatting (bold font, border, centered) onto it so it
# matches the rest of the header row (style index 1 in the sheet).
$headerCell = $ws.Range("U1")
$headerCell.NumberFormat = "@"
$headerCell.Value = "2025-11-22"
$ws.Range("T1").Copy()
$headerCell.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New per-student cells for the 2025-11-22 column: everyone is absent.
$ws.Range("U2").Value = "❌"
$ws.Range("U3").Value = "❌"
$ws.Range("U4").Value = "❌"
$ws.Range("U5").Value = "❌"
$ws.Range("U6").Value = "❌"
$ws.Range("U7").Value = "❌"
$ws.Range("U8").Value = "❌"
$ws.Range("U9").Value = "❌"

# Bump the "Total" column (S) by 1 for each student row since a new date
# column was added to the tracked range.
$ws.Range("S2").Value = 16
$ws.Range("S3").Value = 16
$ws.Range("S4").Value = 16
$ws.Range("S5").Value = 16
$ws.Range("S6").Value = 16
$ws.Range("S7").Value = 16
$ws.Range("S8").Value = 1
$ws.Range("S9").Value = 16
